$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G, shifting TravellerName..EmailID right by one.
$ws.Columns("G").Insert()

# Set explicit column width for the new column (matches F's width, no bestFit).
$ws.Columns("G").ColumnWidth = 8.83

# Header for the new column
$ws.Range("G1").Value = "ToDate"

# Data values for the new column - stored as text (like the other date columns) to avoid
# Excel auto-converting them into date serials.
$ws.Range("G2:G5").NumberFormat = "@"
$ws.Range("G2").Value = "6/12/2017"
$ws.Range("G3").Value = "20/12/2017"
$ws.Range("G4").Value = "20/12/2017"
$ws.Range("G5").Value = "20/12/2017"

$ws.Range("G5").Select()
